$wb = $excel.ActiveWorkbook

$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/f2bc8ee3af7a6e644f3af8b776bb464d5de1467e/e2e/76c684fe-fae5-4582-bbe6-b7188850f6f9.md"
$targetDisplay = "76c684fe-fae5-4582-bbe6-b7188850f6f9.md"

# ---------------- zh-cn sheet ----------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Add($wsZh.Range("J5"), $targetUrl, [Type]::Missing, [Type]::Missing, $targetDisplay) | Out-Null
$wsZh.Range("J5").Style = "Hyperlink"

$wsZh.Range("K5").Value = "76c684fe-fae5-4582-bbe6-b7188850f6f9.f42cdad4abbdd7b81b68f71f93b29284f3e58e70.zh-cn.xlf"
$wsZh.Range("L5").Value = "2017-02-21 02:44:53"
$wsZh.Range("R5").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/2fe4e2f21487c855ee715a7b44834df25c3046a2/e2e/76c684fe-fae5-4582-bbe6-b7188850f6f9.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/f2bc8ee3af7a6e644f3af8b776bb464d5de1467e/e2e/76c684fe-fae5-4582-bbe6-b7188850f6f9.md."

$wsZh.Columns.Item(10).ColumnWidth = 39.17
$wsZh.Columns.Item(11).ColumnWidth = 39.17
$wsZh.Columns.Item(18).ColumnWidth = 39.17

# ---------------- de-de sheet ----------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add($wsDe.Range("J5"), $targetUrl, [Type]::Missing, [Type]::Missing, $targetDisplay) | Out-Null
$wsDe.Range("J5").Style = "Hyperlink"

$wsDe.Range("K5").Value = "76c684fe-fae5-4582-bbe6-b7188850f6f9.f42cdad4abbdd7b81b68f71f93b29284f3e58e70.de-de.xlf"
$wsDe.Range("L5").Value = "2017-02-21 02:45:16"
$wsDe.Range("R5").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/2fe4e2f21487c855ee715a7b44834df25c3046a2/e2e/76c684fe-fae5-4582-bbe6-b7188850f6f9.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/f2bc8ee3af7a6e644f3af8b776bb464d5de1467e/e2e/76c684fe-fae5-4582-bbe6-b7188850f6f9.md."

$wsDe.Columns.Item(10).ColumnWidth = 39.17
$wsDe.Columns.Item(11).ColumnWidth = 39.17
$wsDe.Columns.Item(18).ColumnWidth = 39.17
